# Commit: "changing FALSE to False"
# The I2:I32 column held a boolean formula =FALSE() (displayed via a custom
# "TRUE"/"FALSE" number format). The edit replaces that with the literal
# text value "False" in each cell (text/string content instead of a boolean
# formula result), and updates the active selection to the I column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rng = $ws.Range("I2:I32")

# Make the column hold plain text going forward.
$rng.NumberFormat = "@"

# Replace each boolean-formula cell with the literal text "False".
# The leading apostrophe forces Excel to store it as text rather than
# re-interpreting the word "False" as a boolean TRUE/FALSE value.
for ($r = 2; $r -le 32; $r++) {
    $ws.Cells.Item($r, 9).Value2 = "'False"
}

# Match the author's resulting selection (I2:I32, active cell I2).
$ws.Range("I2:I32").Select() | Out-Null
